# Scheduled runner update: refresh market-board price / profit figures
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for a handful of leves across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Columns: H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#          K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 112 - Making Ends Meet
$ws.Cells.Item(112, 8).Value = 1882.6216
$ws.Cells.Item(112, 10).Value = 2005.0605
$ws.Cells.Item(112, 12).Value = 6015.181500000001
$ws.Cells.Item(112, 14).Value = -8231.181500000001
# Row 137 - Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 1210.381
$ws.Cells.Item(137, 9).Value = 1167.1621
$ws.Cells.Item(137, 10).Value = 1530.2
$ws.Cells.Item(137, 11).Value = 3501.4863
$ws.Cells.Item(137, 12).Value = 4590.6
$ws.Cells.Item(137, 13).Value = -951.4863
$ws.Cells.Item(137, 14).Value = -9690.6
# Row 138 - All-night Crafting
$ws.Cells.Item(138, 8).Value = 2575.349
$ws.Cells.Item(138, 9).Value = 1571.4166
$ws.Cells.Item(138, 10).Value = 3193.1538
$ws.Cells.Item(138, 11).Value = 4714.2498
$ws.Cells.Item(138, 12).Value = 9579.4614
$ws.Cells.Item(138, 13).Value = 425.7502000000004
$ws.Cells.Item(138, 14).Value = -19859.4614

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 61 - Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 2285.9167
$ws.Cells.Item(61, 9).Value = 2115.4119
$ws.Cells.Item(61, 10).Value = 2700
$ws.Cells.Item(61, 11).Value = 2115.4119
$ws.Cells.Item(61, 12).Value = 2700
$ws.Cells.Item(61, 13).Value = -1903.4119
$ws.Cells.Item(61, 14).Value = -3124
# Row 64 - Don't Scuttle with Scuta (also adds LeveProfitHQ)
$ws.Cells.Item(64, 8).Value = 30000
$ws.Cells.Item(64, 10).Value = 31000
$ws.Cells.Item(64, 12).Value = 31000
$ws.Cells.Item(64, 14).Value = -31496
# Row 67 - Shielded by Bureaucracy (L) (also adds LeveProfitHQ)
$ws.Cells.Item(67, 8).Value = 30000
$ws.Cells.Item(67, 10).Value = 31000
$ws.Cells.Item(67, 12).Value = 31000
$ws.Cells.Item(67, 14).Value = -32716
# Row 74 - As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 958.0526
$ws.Cells.Item(74, 9).Value = 606.0833
$ws.Cells.Item(74, 10).Value = 1561.4286
$ws.Cells.Item(74, 11).Value = 606.0833
$ws.Cells.Item(74, 12).Value = 1561.4286
$ws.Cells.Item(74, 13).Value = 267.9167
$ws.Cells.Item(74, 14).Value = -3309.4286
# Row 77 - Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 958.0526
$ws.Cells.Item(77, 9).Value = 606.0833
$ws.Cells.Item(77, 10).Value = 1561.4286
$ws.Cells.Item(77, 11).Value = 3030.4165
$ws.Cells.Item(77, 12).Value = 7807.143
$ws.Cells.Item(77, 13).Value = 1337.5835
$ws.Cells.Item(77, 14).Value = -16543.143
# Row 132 - Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 3783.7737
$ws.Cells.Item(132, 9).Value = 4250.3887
$ws.Cells.Item(132, 11).Value = 12751.1661
$ws.Cells.Item(132, 13).Value = -10221.1661
# Row 136 - Metal with Mettle
$ws.Cells.Item(136, 8).Value = 2285.9167
$ws.Cells.Item(136, 9).Value = 2115.4119
$ws.Cells.Item(136, 10).Value = 2700
$ws.Cells.Item(136, 11).Value = 6346.2357
$ws.Cells.Item(136, 12).Value = 8100
$ws.Cells.Item(136, 13).Value = -3796.2357
$ws.Cells.Item(136, 14).Value = -13200

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Row 62 - Barring the Gates to Foundation
$ws.Cells.Item(62, 8).Value = 30713.715
$ws.Cells.Item(62, 10).Value = 30713.715
$ws.Cells.Item(62, 12).Value = 30713.715
$ws.Cells.Item(62, 14).Value = -32085.715
# Row 65 - Starting Young (L)
$ws.Cells.Item(65, 8).Value = 30713.715
$ws.Cells.Item(65, 10).Value = 30713.715
$ws.Cells.Item(65, 12).Value = 92141.145
$ws.Cells.Item(65, 14).Value = -99005.145
# Row 105 - Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 7939678
$ws.Cells.Item(105, 9).Value = 8931700
$ws.Cells.Item(105, 10).Value = 3499
$ws.Cells.Item(105, 11).Value = 8931700
$ws.Cells.Item(105, 12).Value = 3499
$ws.Cells.Item(105, 13).Value = -8929953
$ws.Cells.Item(105, 14).Value = -6993
# Row 134 - Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 2451.7104
$ws.Cells.Item(134, 9).Value = 2058.9614
$ws.Cells.Item(134, 11).Value = 6176.8842
$ws.Cells.Item(134, 13).Value = -3641.8842

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 31 - Wall Not Found
$ws.Cells.Item(31, 8).Value = 1956.3334
$ws.Cells.Item(31, 9).Value = 1504.7142
$ws.Cells.Item(31, 10).Value = 3537
$ws.Cells.Item(31, 11).Value = 1504.7142
$ws.Cells.Item(31, 12).Value = 3537
$ws.Cells.Item(31, 13).Value = -1209.7142
$ws.Cells.Item(31, 14).Value = -4127
# Row 34 - Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 1956.3334
$ws.Cells.Item(34, 9).Value = 1504.7142
$ws.Cells.Item(34, 10).Value = 3537
$ws.Cells.Item(34, 11).Value = 1504.7142
$ws.Cells.Item(34, 12).Value = 3537
$ws.Cells.Item(34, 13).Value = -1302.7142
$ws.Cells.Item(34, 14).Value = -3941
# Row 58 - You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 927856.5
$ws.Cells.Item(58, 9).Value = 1545019
$ws.Cells.Item(58, 10).Value = 2112.75
$ws.Cells.Item(58, 11).Value = 1545019
$ws.Cells.Item(58, 12).Value = 2112.75
$ws.Cells.Item(58, 13).Value = -1544816
$ws.Cells.Item(58, 14).Value = -2518.75
# Row 132 - Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 348200.53
$ws.Cells.Item(132, 9).Value = 398765.38
$ws.Cells.Item(132, 10).Value = 4359.6
$ws.Cells.Item(132, 11).Value = 1196296.14
$ws.Cells.Item(132, 12).Value = 13078.8
$ws.Cells.Item(132, 13).Value = -1193766.14
$ws.Cells.Item(132, 14).Value = -18138.8
# Row 134 - Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 2155.825
$ws.Cells.Item(134, 9).Value = 1547.7858
$ws.Cells.Item(134, 10).Value = 3574.5833
$ws.Cells.Item(134, 11).Value = 4643.357400000001
$ws.Cells.Item(134, 12).Value = 10723.7499
$ws.Cells.Item(134, 13).Value = -2108.357400000001
$ws.Cells.Item(134, 14).Value = -15793.7499
# Row 136 - Turali Quality
$ws.Cells.Item(136, 8).Value = 927856.5
$ws.Cells.Item(136, 9).Value = 1545019
$ws.Cells.Item(136, 10).Value = 2112.75
$ws.Cells.Item(136, 11).Value = 4635057
$ws.Cells.Item(136, 12).Value = 6338.25
$ws.Cells.Item(136, 13).Value = -4632507
$ws.Cells.Item(136, 14).Value = -11438.25
# Row 141 - No Greater Treasure
$ws.Cells.Item(141, 8).Value = 29368.111
$ws.Cells.Item(141, 9).Value = 15000
$ws.Cells.Item(141, 10).Value = 31164.125
$ws.Cells.Item(141, 11).Value = 15000
$ws.Cells.Item(141, 12).Value = 31164.125
$ws.Cells.Item(141, 13).Value = -9820
$ws.Cells.Item(141, 14).Value = -41524.125

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Row 5 - What a Sap
$ws.Cells.Item(5, 8).Value = 1780.5
$ws.Cells.Item(5, 9).Value = 1855.9231
$ws.Cells.Item(5, 10).Value = 800
$ws.Cells.Item(5, 11).Value = 5567.7693
$ws.Cells.Item(5, 12).Value = 2400
$ws.Cells.Item(5, 13).Value = -5455.7693
$ws.Cells.Item(5, 14).Value = -2624
# Row 131 - The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 1072.1461
$ws.Cells.Item(131, 9).Value = 356
$ws.Cells.Item(131, 10).Value = 1217.3108
$ws.Cells.Item(131, 11).Value = 1068
$ws.Cells.Item(131, 12).Value = 3651.9324
$ws.Cells.Item(131, 13).Value = 3972
$ws.Cells.Item(131, 14).Value = -13731.9324
# Row 135 - Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 1780.5
$ws.Cells.Item(135, 9).Value = 1855.9231
$ws.Cells.Item(135, 10).Value = 800
$ws.Cells.Item(135, 11).Value = 16703.3079
$ws.Cells.Item(135, 12).Value = 7200
$ws.Cells.Item(135, 13).Value = -14168.3079
$ws.Cells.Item(135, 14).Value = -12270

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 132 - On Board for Lar
$ws.Cells.Item(132, 8).Value = 1854.9767
$ws.Cells.Item(132, 9).Value = 1255.4062
$ws.Cells.Item(132, 10).Value = 3599.182
$ws.Cells.Item(132, 11).Value = 3766.2186
$ws.Cells.Item(132, 12).Value = 10797.546
$ws.Cells.Item(132, 13).Value = -1236.2186
$ws.Cells.Item(132, 14).Value = -15857.546

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 132 - Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 7806.5
$ws.Cells.Item(132, 9).Value = 10400
$ws.Cells.Item(132, 10).Value = 5213
$ws.Cells.Item(132, 11).Value = 31200
$ws.Cells.Item(132, 12).Value = 15639
$ws.Cells.Item(132, 13).Value = -28670
$ws.Cells.Item(132, 14).Value = -20699
# Row 136 - Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 31564766
$ws.Cells.Item(136, 9).Value = 41668104
$ws.Cells.Item(136, 10).Value = 1254750.6
$ws.Cells.Item(136, 11).Value = 125004312
$ws.Cells.Item(136, 12).Value = 3764251.8
$ws.Cells.Item(136, 13).Value = -125001762
$ws.Cells.Item(136, 14).Value = -3769351.8

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Row 132 - Comfy Cabins
$ws.Cells.Item(132, 8).Value = 4760.2
$ws.Cells.Item(132, 9).Value = 5176
$ws.Cells.Item(132, 11).Value = 15528
$ws.Cells.Item(132, 13).Value = -12998
# Row 136 - Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 1877.1
$ws.Cells.Item(136, 9).Value = 1692.2084
$ws.Cells.Item(136, 11).Value = 5076.6252
$ws.Cells.Item(136, 13).Value = -2526.6252
